$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: rotate the F:V contents of rows 57, 58 and 59 ---
# New row57(F:V) = old row59(F:V)
# New row58(F:V) = old row57(F:V)
# New row59(F:V) = old row58(F:V)
$row57 = $ws.Range("F57:V57").Value2
$row58 = $ws.Range("F58:V58").Value2
$row59 = $ws.Range("F59:V59").Value2

$ws.Range("F57:V57").Value2 = $row59
$ws.Range("F58:V58").Value2 = $row57
$ws.Range("F59:V59").Value2 = $row58

# --- Part 2: append a new match row (row 97) at the end of the table ---
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Duplicate the formatting (styles) of the previous last row onto the new row
$ws.Range("A" + $lastRow + ":V" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":V" + $newRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A" + $newRow).Value2 = ($newRow - 1)
$ws.Range("B" + $newRow).Value2 = "croatia"
$ws.Range("C" + $newRow).Value2 = "prva-nl"
$ws.Range("D" + $newRow).Value2 = "2023-2024"
$ws.Range("E" + $newRow).Value2 = 45257.70833333334
$ws.Range("F" + $newRow).Value2 = "Sibenik"
$ws.Range("G" + $newRow).Value2 = 1
$ws.Range("H" + $newRow).Value2 = "Orijent"
$ws.Range("I" + $newRow).Value2 = 0
$ws.Range("J" + $newRow).Value2 = 1.23
$ws.Range("K" + $newRow).Value2 = "24/11/2023 01:42"
$ws.Range("L" + $newRow).Value2 = 1.17
$ws.Range("M" + $newRow).Value2 = "27/11/2023 16:58"
$ws.Range("N" + $newRow).Value2 = 5.61
$ws.Range("O" + $newRow).Value2 = "24/11/2023 01:42"
$ws.Range("P" + $newRow).Value2 = 6.98
$ws.Range("Q" + $newRow).Value2 = "27/11/2023 16:58"
$ws.Range("R" + $newRow).Value2 = 8.050000000000001
$ws.Range("S" + $newRow).Value2 = "24/11/2023 01:42"
$ws.Range("T" + $newRow).Value2 = 13.49
$ws.Range("U" + $newRow).Value2 = "27/11/2023 16:58"
$ws.Range("V" + $newRow).Value2 = "https://www.betexplorer.com/football/croatia/prva-nl/sibenik-orijent/pbZmqg1l/"
